$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('E2').Value = '2026-02-06 14:17:43'
$ws.Range('K2').Value = '7.3 MJ/m2'
$ws.Range('M2').Value = '5.3 °C 13:59 TU'
$ws.Range('O2').Value = '-0.6 °C'
$ws.Range('E3').Value = '2026-02-06 14:17:45'
$ws.Range('H3').Value = '71%'
$ws.Range('K3').Value = '9.9 MJ/m2'
$ws.Range('O3').Value = '-2.1 °C'
$ws.Range('E4').Value = '2026-02-06 14:17:48'
$ws.Range('J4').Value = '996.4 hPa'
$ws.Range('K4').Value = '9.8 MJ/m2'
$ws.Range('M4').Value = '18.5 °C 13:47 TU'
$ws.Range('O4').Value = '13.1 °C'
$ws.Range('E5').Value = '2026-02-06 14:17:50'
$ws.Range('H5').Value = '70%'
$ws.Range('J5').Value = '996.7 hPa'
$ws.Range('K5').Value = '9.2 MJ/m2'
$ws.Range('M5').Value = '16.8 °C 13:32 TU'
$ws.Range('O5').Value = '10.3 °C'
$ws.Range('E6').Value = '2026-02-06 14:17:53'
$ws.Range('H6').Value = '47%'
$ws.Range('J6').Value = '997.9 hPa'
$ws.Range('K6').Value = '8.3 MJ/m2'
$ws.Range('O6').Value = '15.2 °C'
$ws.Range('E7').Value = '2026-02-06 14:17:55'
$ws.Range('H7').Value = '63%'
$ws.Range('J7').Value = '997.5 hPa'
$ws.Range('K7').Value = '10.4 MJ/m2'
$ws.Range('M7').Value = '17.0 °C 13:47 TU'
$ws.Range('O7').Value = '11.2 °C'
$ws.Range('E8').Value = '2026-02-06 14:17:57'
$ws.Range('H8').Value = '79%'
$ws.Range('K8').Value = '10.1 MJ/m2'
$ws.Range('O8').Value = '9.4 °C'
$ws.Range('E9').Value = '2026-02-06 14:18:00'
$ws.Range('H9').Value = '89%'
$ws.Range('M9').Value = '12.6 °C 13:55 TU'
$ws.Range('O9').Value = '3.8 °C'
$ws.Range('E10').Value = '2026-02-06 14:18:02'
$ws.Range('H10').Value = '90%'
$ws.Range('M10').Value = '16.1 °C 13:31 TU'
$ws.Range('O10').Value = '8.1 °C'
$ws.Range('E11').Value = '2026-02-06 14:18:04'
$ws.Range('H11').Value = '80%'
$ws.Range('K11').Value = '7.0 MJ/m2'
$ws.Range('O11').Value = '4.7 °C'
$ws.Range('E12').Value = '2026-02-06 14:18:07'
$ws.Range('E13').Value = '2026-02-06 14:18:09'
$ws.Range('H13').Value = '79%'
$ws.Range('O13').Value = '9.4 °C'
$ws.Range('E14').Value = '2026-02-06 14:18:11'
$ws.Range('K14').Value = '5.4 MJ/m2'
$ws.Range('E15').Value = '2026-02-06 14:18:14'
$ws.Range('H15').Value = '75%'
$ws.Range('J15').Value = '996.8 hPa'
$ws.Range('K15').Value = '9.7 MJ/m2'
$ws.Range('M15').Value = '18.5 °C 13:56 TU'
$ws.Range('O15').Value = '9.5 °C'
$ws.Range('E16').Value = '2026-02-06 14:18:16'
$ws.Range('H16').Value = '88%'
$ws.Range('K16').Value = '8.2 MJ/m2'
$ws.Range('O16').Value = '5.4 °C'
$ws.Range('E17').Value = '2026-02-06 14:18:19'
$ws.Range('H17').Value = '89%'
$ws.Range('K17').Value = '8.5 MJ/m2'
$ws.Range('O17').Value = '5.1 °C'
$ws.Range('E18').Value = '2026-02-06 14:18:21'
$ws.Range('K18').Value = '4.8 MJ/m2'
$ws.Range('O18').Value = '-4.4 °C'
$ws.Range('E19').Value = '2026-02-06 14:18:23'
$ws.Range('H19').Value = '80%'
$ws.Range('J19').Value = '999.4 hPa'
$ws.Range('K19').Value = '9.5 MJ/m2'
$ws.Range('M19').Value = '15.5 °C 13:42 TU'
$ws.Range('O19').Value = '9.1 °C'
$ws.Range('E20').Value = '2026-02-06 14:18:25'
$ws.Range('K20').Value = '9.4 MJ/m2'
$ws.Range('O20').Value = '-2.0 °C'
$ws.Range('E21').Value = '2026-02-06 14:18:28'
$ws.Range('H21').Value = '76%'
$ws.Range('J21').Value = '997.2 hPa'
$ws.Range('K21').Value = '9.1 MJ/m2'
$ws.Range('O21').Value = '7.6 °C'
$ws.Range('E22').Value = '2026-02-06 14:18:30'
$ws.Range('H22').Value = '77%'
$ws.Range('K22').Value = '9.3 MJ/m2'
$ws.Range('O22').Value = '10.0 °C'
$ws.Range('E23').Value = '2026-02-06 14:18:33'
$ws.Range('H23').Value = '82%'
$ws.Range('J23').Value = '996.8 hPa'
$ws.Range('K23').Value = '8.3 MJ/m2'
$ws.Range('O23').Value = '9.6 °C'
$ws.Range('E24').Value = '2026-02-06 14:18:35'
$ws.Range('J24').Value = '996.2 hPa'
$ws.Range('K24').Value = '9.8 MJ/m2'
$ws.Range('M24').Value = '16.1 °C 13:59 TU'
$ws.Range('O24').Value = '12.9 °C'
$ws.Range('E25').Value = '2026-02-06 14:18:37'
$ws.Range('H25').Value = '83%'
$ws.Range('I25').Value = '0.1 mm'
$ws.Range('K25').Value = '6.9 MJ/m2'
$ws.Range('O25').Value = '3.6 °C'
$ws.Range('E26').Value = '2026-02-06 14:18:40'
$ws.Range('H26').Value = '77%'
$ws.Range('K26').Value = '7.1 MJ/m2'
$ws.Range('O26').Value = '-1.0 °C'
$ws.Range('E27').Value = '2026-02-06 14:18:42'
$ws.Range('H27').Value = '86%'
$ws.Range('J27').Value = '996.8 hPa'
$ws.Range('K27').Value = '8.7 MJ/m2'
$ws.Range('L27').Value = '19.1 km/h - 104º 13:37 TU'
$ws.Range('O27').Value = '10.0 °C'
$ws.Range('E28').Value = '2026-02-06 14:18:44'
$ws.Range('O28').Value = '4.0 °C'
$ws.Range('E29').Value = '2026-02-06 14:18:47'
$ws.Range('H29').Value = '62%'
$ws.Range('K29').Value = '10.4 MJ/m2'
$ws.Range('O29').Value = '12.2 °C'
$ws.Range('E30').Value = '2026-02-06 14:18:49'
$ws.Range('K30').Value = '7.8 MJ/m2'
$ws.Range('E31').Value = '2026-02-06 14:18:52'
$ws.Range('H31').Value = '89%'
$ws.Range('M31').Value = '13.3 °C 13:43 TU'
$ws.Range('O31').Value = '6.5 °C'
$ws.Range('E32').Value = '2026-02-06 14:18:54'
$ws.Range('J32').Value = '998.3 hPa'
$ws.Range('K32').Value = '10.0 MJ/m2'
$ws.Range('E33').Value = '2026-02-06 14:18:56'
$ws.Range('H33').Value = '85%'
$ws.Range('O33').Value = '9.4 °C'
$ws.Range('E34').Value = '2026-02-06 14:18:59'
$ws.Range('H34').Value = '78%'
$ws.Range('K34').Value = '9.7 MJ/m2'
$ws.Range('O34').Value = '8.0 °C'
$ws.Range('E35').Value = '2026-02-06 14:19:01'
$ws.Range('K35').Value = '7.2 MJ/m2'
$ws.Range('E36').Value = '2026-02-06 14:19:03'
$ws.Range('J36').Value = '999.4 hPa'
$ws.Range('K36').Value = '9.8 MJ/m2'
$ws.Range('O36').Value = '12.8 °C'
